$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 657.4706
$ws.Range("I41").Value = 607.0909
$ws.Range("K41").Value = 607.0909
$ws.Range("M41").Value = -167.0909
$ws.Range("H51").Value = 2588.5
$ws.Range("I51").Value = 2200
$ws.Range("K51").Value = 2200
$ws.Range("M51").Value = -1716
$ws.Range("H70").Value = 58816.723
$ws.Range("I70").Value = 2350.5
$ws.Range("J70").Value = 65875
$ws.Range("K70").Value = 7051.5
$ws.Range("L70").Value = 197625
$ws.Range("M70").Value = -6781.5
$ws.Range("N70").Value = -198165
$ws.Range("H73").Value = 58816.723
$ws.Range("I73").Value = 2350.5
$ws.Range("J73").Value = 65875
$ws.Range("K73").Value = 7051.5
$ws.Range("L73").Value = 197625
$ws.Range("M73").Value = -6115.5
$ws.Range("N73").Value = -199497
$ws.Range("H76").Value = 4604.8125
$ws.Range("I76").Value = 3799.889
$ws.Range("K76").Value = 3799.889
$ws.Range("M76").Value = -3484.889
$ws.Range("H79").Value = 4604.8125
$ws.Range("I79").Value = 3799.889
$ws.Range("K79").Value = 3799.889
$ws.Range("M79").Value = -2707.889
$ws.Range("H92").Value = 1295.0625
$ws.Range("I92").Value = 1295.0625
$ws.Range("K92").Value = 1295.0625
$ws.Range("M92").Value = -47.0625
$ws.Range("H111").Value = 3255.6667
$ws.Range("I111").Value = 581.2857
$ws.Range("J111").Value = 6999.8
$ws.Range("K111").Value = 1743.8571
$ws.Range("L111").Value = 20999.4
$ws.Range("M111").Value = 1323.1429
$ws.Range("N111").Value = -27133.4
$ws.Range("H113").Value = 7552.533
$ws.Range("I113").Value = 5161
$ws.Range("K113").Value = 5161
$ws.Range("M113").Value = -1907
$ws.Range("H116").Value = 17174.234
$ws.Range("I116").Value = 16497.467
$ws.Range("J116").Value = 22250
$ws.Range("K116").Value = 16497.467
$ws.Range("L116").Value = 22250
$ws.Range("M116").Value = -13055.467
$ws.Range("N116").Value = -29134
$ws.Range("H138").Value = 10207615
$ws.Range("J138").Value = 14710420
$ws.Range("L138").Value = 44131260
$ws.Range("N138").Value = -44141540

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3402.0469
$ws.Range("I32").Value = 3402.0469
$ws.Range("K32").Value = 3402.0469
$ws.Range("M32").Value = -3115.0469
$ws.Range("H44").Value = 23249.334
$ws.Range("I44").Value = 14749.5
$ws.Range("K44").Value = 14749.5
$ws.Range("M44").Value = -14261.5
$ws.Range("H60").Value = 10112
$ws.Range("I60").Value = 10112
$ws.Range("K60").Value = 10112
$ws.Range("M60").Value = -9379
$ws.Range("H61").Value = 5581.4375
$ws.Range("I61").Value = 5388.25
$ws.Range("J61").Value = 5645.8335
$ws.Range("K61").Value = 5388.25
$ws.Range("L61").Value = 5645.8335
$ws.Range("M61").Value = -5176.25
$ws.Range("N61").Value = -6069.8335
$ws.Range("H119").Value = 39846
$ws.Range("J119").Value = 39846
$ws.Range("L119").Value = 39846
$ws.Range("N119").Value = -49522
$ws.Range("H136").Value = 5581.4375
$ws.Range("I136").Value = 5388.25
$ws.Range("J136").Value = 5645.8335
$ws.Range("K136").Value = 16164.75
$ws.Range("L136").Value = 16937.5005
$ws.Range("M136").Value = -13614.75
$ws.Range("N136").Value = -22037.5005

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2027.8
$ws.Range("I86").Value = 1932.25
$ws.Range("K86").Value = 1932.25
$ws.Range("M86").Value = -809.25
$ws.Range("H89").Value = 2027.8
$ws.Range("I89").Value = 1932.25
$ws.Range("K89").Value = 9661.25
$ws.Range("M89").Value = -4045.25
$ws.Range("H105").Value = 1577.6111
$ws.Range("I105").Value = 1330.2
$ws.Range("K105").Value = 1330.2
$ws.Range("M105").Value = 416.8
$ws.Range("H107").Value = 9520.154
$ws.Range("I107").Value = 2113.6365
$ws.Range("K107").Value = 2113.6365
$ws.Range("M107").Value = -193.6365000000001
$ws.Range("H135").Value = 56944.5
$ws.Range("J135").Value = 56944.5
$ws.Range("L135").Value = 56944.5
$ws.Range("N135").Value = -67084.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 114231.89
$ws.Range("I31").Value = 202895
$ws.Range("J31").Value = 3403
$ws.Range("K31").Value = 202895
$ws.Range("L31").Value = 3403
$ws.Range("M31").Value = -202600
$ws.Range("N31").Value = -3993
$ws.Range("H34").Value = 114231.89
$ws.Range("I34").Value = 202895
$ws.Range("J34").Value = 3403
$ws.Range("K34").Value = 202895
$ws.Range("L34").Value = 3403
$ws.Range("M34").Value = -202693
$ws.Range("N34").Value = -3807
$ws.Range("H50").Value = 12500
$ws.Range("J50").Value = 12500
$ws.Range("L50").Value = 12500
$ws.Range("N50").Value = -13750
$ws.Range("H134").Value = 27713.56
$ws.Range("I134").Value = 8152.8335
$ws.Range("K134").Value = 24458.5005
$ws.Range("M134").Value = -21923.5005
$ws.Range("H141").Value = 417780
$ws.Range("J141").Value = 417780
$ws.Range("L141").Value = 417780
$ws.Range("N141").Value = -428140

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 985.75
$ws.Range("J107").Value = 1600
$ws.Range("L107").Value = 4800
$ws.Range("N107").Value = -8640
$ws.Range("H133").Value = 15762.25
$ws.Range("I133").Value = 16016.333
$ws.Range("J133").Value = 15000
$ws.Range("K133").Value = 48048.999
$ws.Range("L133").Value = 45000
$ws.Range("M133").Value = -42988.999
$ws.Range("N133").Value = -55120
$ws.Range("H136").Value = 458801
$ws.Range("I136").Value = 589389.5
$ws.Range("K136").Value = 1768168.5
$ws.Range("M136").Value = -1763068.5
$ws.Range("H137").Value = 2510.7144
$ws.Range("I137").Value = 2066
$ws.Range("K137").Value = 6198
$ws.Range("M137").Value = -1098
$ws.Range("H138").Value = 27785930
$ws.Range("I138").Value = 33340116
$ws.Range("K138").Value = 100020348
$ws.Range("M138").Value = -100015208
$ws.Range("H139").Value = 1380.8334
$ws.Range("I139").Value = 1310.4348
$ws.Range("K139").Value = 3931.3044
$ws.Range("M139").Value = 1208.6956
$ws.Range("H140").Value = 2819.8
$ws.Range("I140").Value = 2561.3076
$ws.Range("K140").Value = 7683.9228
$ws.Range("M140").Value = -2503.9228

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9221.75
$ws.Range("J92").Value = 8962.333000000001
$ws.Range("L92").Value = 8962.333000000001
$ws.Range("N92").Value = -12706.333
$ws.Range("H107").Value = 545.9167
$ws.Range("I107").Value = 438.55554
$ws.Range("J107").Value = 868
$ws.Range("K107").Value = 438.55554
$ws.Range("L107").Value = 868
$ws.Range("M107").Value = 1481.44446
$ws.Range("N107").Value = -4708
$ws.Range("H132").Value = 4632.4546
$ws.Range("I132").Value = 4254.3887
$ws.Range("K132").Value = 12763.1661
$ws.Range("M132").Value = -10233.1661

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2014.4
$ws.Range("I22").Value = 2775
$ws.Range("J22").Value = 1507.3334
$ws.Range("K22").Value = 2775
$ws.Range("L22").Value = 1507.3334
$ws.Range("M22").Value = -2480
$ws.Range("N22").Value = -2097.3334
$ws.Range("H27").Value = 2014.4
$ws.Range("I27").Value = 2775
$ws.Range("J27").Value = 1507.3334
$ws.Range("K27").Value = 2775
$ws.Range("L27").Value = 1507.3334
$ws.Range("M27").Value = -2668
$ws.Range("N27").Value = -1721.3334
$ws.Range("H68").Value = 3716.5
$ws.Range("I68").Value = 3533.1428
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 3533.1428
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -2784.1428
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 3716.5
$ws.Range("I71").Value = 3533.1428
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 17665.714
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -13921.714
$ws.Range("N71").Value = -32488
$ws.Range("H122").Value = 630449.0600000001
$ws.Range("I122").Value = 913399.2
$ws.Range("J122").Value = 7958.8
$ws.Range("K122").Value = 2740197.6
$ws.Range("L122").Value = 23876.4
$ws.Range("M122").Value = -2737747.6
$ws.Range("N122").Value = -28776.4
$ws.Range("H133").Value = 60481.5
$ws.Range("J133").Value = 60481.5
$ws.Range("L133").Value = 60481.5
$ws.Range("N133").Value = -65541.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 273.33334
$ws.Range("I107").Value = 290.53845
$ws.Range("J107").Value = 228.6
$ws.Range("K107").Value = 871.61535
$ws.Range("L107").Value = 685.8
$ws.Range("M107").Value = 1048.38465
$ws.Range("N107").Value = -4525.8
$ws.Range("H122").Value = 2394.0952
$ws.Range("I122").Value = 2120.5293
$ws.Range("J122").Value = 3556.75
$ws.Range("K122").Value = 6361.5879
$ws.Range("L122").Value = 10670.25
$ws.Range("M122").Value = -3911.5879
$ws.Range("N122").Value = -15570.25
$ws.Range("H132").Value = 3893.625
$ws.Range("I132").Value = 2783.25
$ws.Range("J132").Value = 5004
$ws.Range("K132").Value = 8349.75
$ws.Range("L132").Value = 15012
$ws.Range("M132").Value = -5819.75
$ws.Range("N132").Value = -20072
